$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Item3Chance" column (L), shifting
# Item3Chance/Item3MinQ/Item3MaxQ one column to the right (L->M, M->N, N->O),
# so the new column L can hold "Item3Name".
$ws.Columns.Item(12).Insert()

# New header + data for the inserted Item3Name column.
$ws.Range("L1").Value() = "Item3Name"
$ws.Range("L2").Value() = "Iron Cuirass"
$ws.Range("L3").Value() = "Iron Cuirass"
$ws.Range("L4").Value() = "Iron Boots"
$ws.Range("L5").Value() = "Iron Spear"

# Give the new column the plain (non-bestFit) width used in the target file.
$ws.Columns.Item(12).ColumnWidth = 10.125

# Build the per-row JSON-ish summary formula in column Q.
$ws.Range("Q2").Formula = '=_xlfn.CONCAT("{",CHAR(34),$A$1,CHAR(34),":",CHAR(34),$A2,CHAR(34),",",CHAR(34),$B$1,CHAR(34),":",CHAR(34),$B2,CHAR(34),",",CHAR(34),$C$1,CHAR(34),":",CHAR(34),$C2,CHAR(34),",",CHAR(34),$D$1,CHAR(34),":",CHAR(34),$D2,CHAR(34),",",CHAR(34),$E$1,CHAR(34),":",CHAR(34),$E2,CHAR(34),",",CHAR(34),$F$1,CHAR(34),":",CHAR(34),$F2,CHAR(34),",",CHAR(34),$G$1,CHAR(34),":",CHAR(34),$G2,CHAR(34),",",CHAR(34),$H$1,CHAR(34),":",CHAR(34),$H2,CHAR(34),",",CHAR(34),$I$1,CHAR(34),":",CHAR(34),$I2,CHAR(34),",",CHAR(34),$J$1,CHAR(34),":",CHAR(34),$J2,CHAR(34),",",CHAR(34),$K$1,CHAR(34),":",CHAR(34),$K2,CHAR(34),",",CHAR(34),$L$1,CHAR(34),":",CHAR(34),$L2,CHAR(34),",",CHAR(34),$M$1,CHAR(34),":",CHAR(34),$M2,CHAR(34),",",CHAR(34),$N$1,CHAR(34),":",CHAR(34),$N2,CHAR(34),",",CHAR(34),$O$1,CHAR(34),":",CHAR(34),$O2,CHAR(34),"},")'

$ws.Range("Q3:Q5").Formula = '=_xlfn.CONCAT("{",CHAR(34),$A$1,CHAR(34),":",CHAR(34),$A3,CHAR(34),",",CHAR(34),$B$1,CHAR(34),":",CHAR(34),$B3,CHAR(34),",",CHAR(34),$C$1,CHAR(34),":",CHAR(34),$C3,CHAR(34),",",CHAR(34),$D$1,CHAR(34),":",CHAR(34),$D3,CHAR(34),",",CHAR(34),$E$1,CHAR(34),":",CHAR(34),$E3,CHAR(34),",",CHAR(34),$F$1,CHAR(34),":",CHAR(34),$F3,CHAR(34),",",CHAR(34),$G$1,CHAR(34),":",CHAR(34),$G3,CHAR(34),",",CHAR(34),$H$1,CHAR(34),":",CHAR(34),$H3,CHAR(34),",",CHAR(34),$I$1,CHAR(34),":",CHAR(34),$I3,CHAR(34),",",CHAR(34),$J$1,CHAR(34),":",CHAR(34),$J3,CHAR(34),",",CHAR(34),$K$1,CHAR(34),":",CHAR(34),$K3,CHAR(34),",",CHAR(34),$L$1,CHAR(34),":",CHAR(34),$L3,CHAR(34),",",CHAR(34),$M$1,CHAR(34),":",CHAR(34),$M3,CHAR(34),",",CHAR(34),$N$1,CHAR(34),":",CHAR(34),$N3,CHAR(34),",",CHAR(34),$O$1,CHAR(34),":",CHAR(34),$O3,CHAR(34),"},")'

# Page setup tweak that accompanied this change.
$ws.PageSetup.Orientation = 1

# Selection / view state to match the target file.
$ws.Range("Q2:Q5").Select()

$wb.Save()
